$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update numeric values in column R (delay/check column) ---
$rUpdates = @{
    2  = -0.3623372395833334
    3  = -0.4377821180555556
    4  = -0.4741319444444445
    6  = 0
    7  = -0.5925130208333333
    8  = -1.416059027777778
    9  = -1.443691314548611
    10 = -1.406147540983796
    11 = 0
    12 = -1.498667800451389
    13 = -0.4885587431712963
    14 = -1.514301215277778
    15 = -0.5466647104861111
    16 = -2.611839708564815
    17 = -2.471889671365741
    18 = -0.5261067708333333
    19 = -2.5365234375
    20 = 0
    22 = -0.4133138020833333
    23 = -1.526104797974537
    24 = -1.442361111111111
    25 = -1.466666666666667
}

foreach ($row in $rUpdates.Keys) {
    $ws.Range("R$row").Value = $rUpdates[$row]
}

# --- Update column N: numeric "veicolo" codes become annotated text ---
$nUpdates = @{
    5  = "39666 (non in estrazione)"
    6  = "39742 (esterno)"
    11 = "39666 (esterno)"
    20 = "39762 (esterno)"
    21 = "39723 (non in estrazione)"
    26 = "39750 (non in estrazione)"
    27 = "39764 (non in estrazione)"
}

foreach ($row in $nUpdates.Keys) {
    $ws.Range("N$row").Value = $nUpdates[$row]
}

# --- Remove the "NESSUN VEICOLO (...)" text cells in column N for rows 31-90 ---
for ($row = 31; $row -le 90; $row++) {
    $ws.Range("N$row").ClearContents()
}
